$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Jan 30 2023 07:10:59 UTC symbol-list refresh from the coinranking.com source feed:
# refreshed Price/Volume(1h) figures, re-ranked several exchange-token rows (their
# Coin/Link moved up or down one slot), and advanced the Hora column from 6 to 7.
# D/E/G are stored as text in the sheet, so values are written with a leading
# apostrophe to keep Excel from re-interpreting numeric-looking strings as numbers.

# Row 2
$ws.Range("D2").Value = '''314.58'
$ws.Range("E2").Value = '''1.95%'
$ws.Range("G2").Value = '''7'

# Row 3
$ws.Range("D3").Value = '''39.25'
$ws.Range("E3").Value = '''-0.72%'
$ws.Range("G3").Value = '''7'

# Row 4
$ws.Range("D4").Value = '''5.139'
$ws.Range("E4").Value = '''0.08%'
$ws.Range("G4").Value = '''7'

# Row 5
$ws.Range("D5").Value = '''0.08164'
$ws.Range("E5").Value = '''0.46%'
$ws.Range("G5").Value = '''7'

# Row 6
$ws.Range("D6").Value = '''1.969'
$ws.Range("E6").Value = '''1.03%'
$ws.Range("G6").Value = '''7'

# Row 7
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").Value = '''8.319'
$ws.Range("E7").Value = '''2.25%'
$ws.Range("G7").Value = '''7'

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.9399'
$ws.Range("E8").Value = '''1.45%'
$ws.Range("G8").Value = '''7'

# Row 9
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '''0.1315'
$ws.Range("E9").Value = '''-7.00%'
$ws.Range("G9").Value = '''7'

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1975'
$ws.Range("E10").Value = '''2.22%'
$ws.Range("G10").Value = '''7'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.09056'
$ws.Range("E11").Value = '''-0.45%'
$ws.Range("G11").Value = '''7'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.03489'
$ws.Range("E12").Value = '''-0.70%'
$ws.Range("G12").Value = '''7'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09715'
$ws.Range("E13").Value = '''-1.13%'
$ws.Range("G13").Value = '''7'

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001416'
$ws.Range("E14").Value = '''1.74%'
$ws.Range("G14").Value = '''7'

# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.006012'
$ws.Range("E15").Value = '''2.56%'
$ws.Range("G15").Value = '''7'

# Row 16
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.608'
$ws.Range("E16").Value = '''-7.69%'
$ws.Range("G16").Value = '''7'

# Row 17
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''4.449'
$ws.Range("E17").Value = '''5.09%'
$ws.Range("G17").Value = '''7'

# Row 18
$ws.Range("D18").Value = '''3.189'
$ws.Range("E18").Value = '''-4.65%'
$ws.Range("G18").Value = '''7'

# Row 19
$ws.Range("D19").Value = '''0.3468'
$ws.Range("E19").Value = '''0.40%'
$ws.Range("G19").Value = '''7'

# Row 20
$ws.Range("E20").Value = '''-0.25%'
$ws.Range("G20").Value = '''7'

# Row 21
$ws.Range("D21").Value = '''5.013'
$ws.Range("E21").Value = '''5.80%'
$ws.Range("G21").Value = '''7'

# Row 22
$ws.Range("D22").Value = '''0.2491'
$ws.Range("E22").Value = '''2.67%'
$ws.Range("G22").Value = '''7'

# Row 23
$ws.Range("D23").Value = '''0.04375'
$ws.Range("E23").Value = '''-0.07%'
$ws.Range("G23").Value = '''7'

# Row 24
$ws.Range("D24").Value = '''0.001243'
$ws.Range("E24").Value = '''1.06%'
$ws.Range("G24").Value = '''7'

# Row 25
$ws.Range("D25").Value = '''0.004748'
$ws.Range("G25").Value = '''7'

# Row 26
$ws.Range("E26").Value = '''199.33%'
$ws.Range("G26").Value = '''7'

# Row 27
$ws.Range("E27").Value = '''-7.62%'
$ws.Range("G27").Value = '''7'

# Row 28
$ws.Range("G28").Value = '''7'

# Row 29
$ws.Range("G29").Value = '''7'

# Row 30
$ws.Range("G30").Value = '''7'

# Row 31
$ws.Range("G31").Value = '''7'

# Row 32
$ws.Range("G32").Value = '''7'

# Row 33
$ws.Range("G33").Value = '''7'

# Row 34
$ws.Range("G34").Value = '''7'

# Row 35
$ws.Range("G35").Value = '''7'

# Row 36
$ws.Range("G36").Value = '''7'

# Row 37
$ws.Range("G37").Value = '''7'

# Row 38
$ws.Range("G38").Value = '''7'

# Row 39
$ws.Range("D39").Value = '''0.02244'
$ws.Range("E39").Value = '''8.13%'
$ws.Range("G39").Value = '''7'

# Row 40
$ws.Range("D40").Value = '''0.05245'
$ws.Range("E40").Value = '''2.89%'
$ws.Range("G40").Value = '''7'

# Row 41
$ws.Range("D41").Value = '''0.007570'
$ws.Range("E41").Value = '''1.59%'
$ws.Range("G41").Value = '''7'

# Row 42
$ws.Range("D42").Value = '''0.01031'
$ws.Range("E42").Value = '''5.32%'
$ws.Range("G42").Value = '''7'

# Row 43
$ws.Range("D43").Value = '''0.1395'
$ws.Range("E43").Value = '''2.19%'
$ws.Range("G43").Value = '''7'

# Row 44
$ws.Range("E44").Value = '''-1.38%'
$ws.Range("G44").Value = '''7'

# Row 45
$ws.Range("D45").Value = '''0.009128'
$ws.Range("E45").Value = '''-4.47%'
$ws.Range("G45").Value = '''7'

# Row 46
$ws.Range("D46").Value = '''0.00006839'
$ws.Range("E46").Value = '''6.75%'
$ws.Range("G46").Value = '''7'

# Row 47
$ws.Range("E47").Value = '''0.04%'
$ws.Range("G47").Value = '''7'

# Row 48
$ws.Range("D48").Value = '''0.003016'
$ws.Range("E48").Value = '''13.55%'
$ws.Range("G48").Value = '''7'

# Row 49
$ws.Range("G49").Value = '''7'

# Row 50
$ws.Range("D50").Value = '''0.00002103'
$ws.Range("E50").Value = '''0.04%'
$ws.Range("G50").Value = '''7'

# Row 51
$ws.Range("D51").Value = '''0.0002003'
$ws.Range("E51").Value = '''0.04%'
$ws.Range("G51").Value = '''7'
